$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the hourly production values (row 4 for Line 1 / Line 2, row 5 H column for Line 2)
$ws.Range("B4").Value = 141
$ws.Range("C4").Value = 154
$ws.Range("G4").Value = 157
$ws.Range("H5").Value = 154

# Move the active selection to I18 (outside the used data range), matching the saved view state
$ws.Range("I18").Select()
